# Update countries & provincias Spain
# - Corrects the ordering of several country names that were in the wrong
#   rows (the per-country stats stay matched to the correct country name).
# - Refreshes the case/death/recovered counters for the affected countries.
# - Bumps the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix country names that were swapped with their neighbour row ---
$ws.Range("A32").Value = "Japon"
$ws.Range("A33").Value = "Bielorrusia"

$ws.Range("A55").Value = "Argentina"
$ws.Range("A56").Value = "Marruecos"

$ws.Range("A124").Value = "Venezuela"
$ws.Range("A125").Value = "Mauricio"

$ws.Range("A158").Value = "Haiti"
$ws.Range("A159").Value = "Bahamas"
$ws.Range("A160").Value = "Barbados"
$ws.Range("A161").Value = "Mozambique"

# --- Refresh the statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) for every row whose
#     numbers changed in this update ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1094730
$ws.Range("C4").Value = 30536
$ws.Range("D4").Value = 152324
$ws.Range("E4").Value = 878565
$ws.Range("F4").Value = 15226
$ws.Range("G4").Value = 2186
$ws.Range("H4").Value = 63841

# Row 32 - Japon
$ws.Range("B32").Value = 14088
$ws.Range("C32").Value = 193
$ws.Range("D32").Value = 2460
$ws.Range("E32").Value = 11198
$ws.Range("F32").Value = 308
$ws.Range("G32").Value = 17
$ws.Range("H32").Value = 430

# Row 33 - Bielorrusia
$ws.Range("B33").Value = 14027
$ws.Range("C33").Value = 846
$ws.Range("D33").Value = 2386
$ws.Range("E33").Value = 11552
$ws.Range("F33").Value = 92
$ws.Range("G33").Value = 5
$ws.Range("H33").Value = 89

# Row 55 - Argentina
$ws.Range("B55").Value = 4428
$ws.Range("C55").Value = 143
$ws.Range("D55").Value = 1256
$ws.Range("E55").Value = 2954
$ws.Range("F55").Value = 157
$ws.Range("G55").Value = 4
$ws.Range("H55").Value = 218

# Row 56 - Marruecos
$ws.Range("B56").Value = 4423
$ws.Range("C56").Value = 102
$ws.Range("D56").Value = 984
$ws.Range("E56").Value = 3269
$ws.Range("F56").Value = 1
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 170

# Row 101 - Niger
$ws.Range("B101").Value = 719
$ws.Range("C101").Value = 6
$ws.Range("D101").Value = 452
$ws.Range("E101").Value = 235

# Row 124 - Venezuela
$ws.Range("B124").Value = 333
$ws.Range("C124").Value = 2
$ws.Range("D124").Value = 142
$ws.Range("E124").Value = 175
$ws.Range("F124").Value = 2
$ws.Range("G124").Value = 6
$ws.Range("H124").Value = 16

# Row 125 - Mauricio
$ws.Range("B125").Value = 332
$ws.Range("D125").Value = 310
$ws.Range("E125").Value = 12
$ws.Range("F125").Value = 3

# Row 158 - Haiti
$ws.Range("C158").Value = 5
$ws.Range("D158").Value = 8
$ws.Range("E158").Value = 65
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 2
$ws.Range("H158").Value = 8

# Row 159 - Bahamas
$ws.Range("B159").Value = 81
$ws.Range("C159").Value = 1
$ws.Range("E159").Value = 45

# Row 160 - Barbados
$ws.Range("B160").Value = 81
$ws.Range("C160").Value = 1
$ws.Range("D160").Value = 39
$ws.Range("E160").Value = 35
$ws.Range("F160").Value = 4
$ws.Range("H160").Value = 7

# Row 161 - Mozambique
$ws.Range("D161").Value = 12
$ws.Range("E161").Value = 64
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 0

# --- Bump the "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 01:52"
